$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (shifts existing rows 2-4 down to 3-5),
# carrying down formatting the same way Excel's own Insert does.
$ws.Rows.Item(2).Insert()

# Populate the newly-inserted row 2 with the new recognition entry.
$ws.Range("A2").Value = "Premio extraordinario de doctorado 2021 - 2022"
$ws.Range("B2").Value = 2023
$ws.Range("C2").Value = "Universidad Carlos III de Madrid "
$ws.Range("D2").Value = "Madrid, España"

# Match the saved selection/active cell recorded in the target workbook.
$ws.Range("A10").Select()
